# skills210322.xlsx — add normalized (no-space/punctuation) tokens for the
# multi-word / symbol-bearing skill labels so they can be used as, e.g., R
# document-term-matrix column names. The sheet keeps column A as the master
# sorted skill list; any skill whose "clean" token differs from its original
# label also gets the original label in column B and a duplicate of the
# clean token in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: introduce the 11 new normalized skill tokens into the shared-string
# table in the same top-to-bottom order the human-readable labels they
# correspond to already occupied (original B42:B52), mirroring how the
# author appears to have generated them (e.g. via a helper column) before
# re-sorting the sheet.
$newTokenOrder = @(
    "csharp",
    "cplusplus",
    "computerscience",
    "dataengineering",
    "deeplearning",
    "machinelearning",
    "neuralnetworks",
    "projectmanagement",
    "scikitlearn",
    "softwaredevelopment",
    "softwareengineering"
)
$tmpRow = 60
foreach ($tok in $newTokenOrder) {
    $ws.Cells.Item($tmpRow, 5).Value = $tok
    $tmpRow = $tmpRow + 1
}

# Target data: row index, column A value (normalized skill token),
# optional column B value (original human-readable label).
# Column C (when present) duplicates column A.
$rows = @(
    @(1, "ai", $null),
    @(2, "analysis", $null),
    @(3, "aws", $null),
    @(4, "azure", $null),
    @(5, "c", $null),
    @(6, "caffe", $null),
    @(7, "cassandra", $null),
    @(8, "communication", $null),
    @(9, "computerscience", "computer science"),
    @(10, "cplusplus", "c++"),
    @(11, "csharp", "c#"),
    @(12, "d3", $null),
    @(13, "dataengineering", "data engineering"),
    @(14, "deeplearning", "deep learning"),
    @(15, "docker", $null),
    @(16, "excel", $null),
    @(17, "git", $null),
    @(18, "hadoop", $null),
    @(19, "hbase", $null),
    @(20, "hive", $null),
    @(21, "java", $null),
    @(22, "javascript", $null),
    @(23, "keras", $null),
    @(24, "linux", $null),
    @(25, "machinelearning", "machine learning"),
    @(26, "mathematics", $null),
    @(27, "matlab", $null),
    @(28, "mongodb", $null),
    @(29, "mysql", $null),
    @(30, "neuralnetworks", "neural networks"),
    @(31, "nlp", $null),
    @(32, "nosql", $null),
    @(33, "numpy", $null),
    @(34, "pandas", $null),
    @(35, "perl", $null),
    @(36, "pig", $null),
    @(37, "projectmanagement", "project management"),
    @(38, "python", $null),
    @(39, "pytorch", $null),
    @(40, "r", $null),
    @(41, "sas", $null),
    @(42, "scala", $null),
    @(43, "scikitlearn", "scikit-learn"),
    @(44, "softwaredevelopment", "software development"),
    @(45, "softwareengineering", "software engineering"),
    @(46, "spark", $null),
    @(47, "spss", $null),
    @(48, "sql", $null),
    @(49, "statistics", $null),
    @(50, "tableau", $null),
    @(51, "tensorflow", $null),
    @(52, "visualization", $null)
)

# Clear the previous A:B data range and the scratch column before laying out
# the final A:C range.
$ws.Range("A1:C60").ClearContents() | Out-Null
$ws.Range("E60:E70").ClearContents() | Out-Null

foreach ($row in $rows) {
    $r = $row[0]
    $aVal = $row[1]
    $bVal = $row[2]

    $ws.Cells.Item($r, 1).Value = $aVal
    if ($bVal -ne $null) {
        $ws.Cells.Item($r, 2).Value = $bVal
        $ws.Cells.Item($r, 3).Value = $aVal
    }
}

# Refresh the sheet's recorded sort state/selection over the new A1:C52 extent.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("A1:A52")) | Out-Null
$sortObj.SetRange($ws.Range("A1:C52"))
$sortObj.Header = 0
$sortObj.Apply()

$sel = $ws.Range("A1:A52")
$sel.Select() | Out-Null
try {
    $ws.Application.ActiveWindow.ScrollRow = 11
} catch {
    # Window scroll position is cosmetic / session-local; ignore if the
    # host doesn't expose it.
}

Write-Output "done"
